$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, pushing existing data (rows 51-114) down
# to rows 52-115, matching the author's weekly price update.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the latest weekly record
# (same dataset template as the former row 51, with updated date/volume/price).
$ws.Range("A51").Value = 6
$ws.Range("B51").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C51").Value = "Metropolitana"
$ws.Range("D51").Value = 45175
$ws.Range("E51").Value = 13
$ws.Range("F51").Value = 100112035
$ws.Range("G51").Value = "Bruselas (repollito)"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 18000
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = 19133
$ws.Range("N51").Value = "`$/malla 15 kilos"
$ws.Range("O51").Value = "Provincia de Quillota"
$ws.Range("P51").Value = 1276
$ws.Range("Q51").Value = 15
$ws.Range("R51").Value = "Hortaliza"
